$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "27.923.91"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "1.764.97"
$ws.Range("E3").Value = "  +0.78%  "

Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue "D5" "329.12"
$ws.Range("E5").Value = "  +1.53%  "

Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -0.04%  "

Set-TextValue "D7" "0.4544"
$ws.Range("E7").Value = "  +1.25%  "

Set-TextValue "D8" "0.3516"
$ws.Range("E8").Value = "  -1.24%  "

Set-TextValue "D9" "42.05"
$ws.Range("E9").Value = "  +1.75%  "

Set-TextValue "D10" "0.07382"
$ws.Range("E10").Value = "  -0.96%  "

Set-TextValue "D11" "1.097"
$ws.Range("E11").Value = "  +1.44%  "

$ws.Range("E12").Value = "  +0.02%  "

Set-TextValue "D13" "20.72"
$ws.Range("E13").Value = "  +0.01%  "

Set-TextValue "D14" "5.995"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("D16").Value = "1.767.60"
$ws.Range("E16").Value = "  +0.94%  "

Set-TextValue "D17" "92.38"
$ws.Range("E17").Value = "  -1.34%  "

$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("E20").Value = "  +0.01%  "

Set-TextValue "D21" "16.95"
$ws.Range("E21").Value = "  -0.61%  "

Set-TextValue "D22" "5.775"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").Value = "27.949.35"
$ws.Range("E23").Value = "  +1.24%  "

Set-TextValue "D24" "11.22"
$ws.Range("E24").Value = "  +0.45%  "

Set-TextValue "D25" "2.151"
$ws.Range("E25").Value = "  +3.16%  "

Set-TextValue "D26" "162.32"

Set-TextValue "D27" "20.15"
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").Value = "1.972.61"
$ws.Range("E28").Value = "  +0.98%  "

Set-TextValue "D29" "2.157"
$ws.Range("E29").Value = "  +2.66%  "

Set-TextValue "D30" "123.84"
$ws.Range("E30").Value = "  -1.03%  "

Set-TextValue "D31" "1.072"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("E32").Value = "  +1.41%  "

Set-TextValue "D33" "3.667"
$ws.Range("E33").Value = "  +0.39%  "

Set-TextValue "D34" "5.583"
$ws.Range("E34").Value = "  +1.70%  "

$ws.Range("E35").Value = "  +0.86%  "

Set-TextValue "D36" "0.02275"
$ws.Range("E36").Value = "  -0.27%  "

Set-TextValue "D37" "0.06121"
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("E38").Value = "  -0.02%  "

Set-TextValue "D39" "4.935"
$ws.Range("E39").Value = "  +0.33%  "

Set-TextValue "D40" "0.6249"
$ws.Range("E40").Value = "  -0.28%  "

Set-TextValue "D41" "1.183"
$ws.Range("E41").Value = "  +0.54%  "

Set-TextValue "D42" "1.382"
$ws.Range("E42").Value = "  -1.05%  "

Set-TextValue "D43" "7.807"
$ws.Range("E43").Value = "  +0.56%  "

Set-TextValue "D44" "13.11"
$ws.Range("E44").Value = "  -0.29%  "

Set-TextValue "D45" "3.735"
$ws.Range("E45").Value = "  +0.54%  "

$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("E48").Value = "  +0.24%  "

Set-TextValue "D49" "1.128"
$ws.Range("E49").Value = "  +0.13%  "

Set-TextValue "D50" "0.06810"
$ws.Range("E50").Value = "  -1.02%  "

Set-TextValue "D51" "72.85"
$ws.Range("E51").Value = "  +2.14%  "
